# excess_mortality_provinces.xlsx — weekly data refresh
# (mirrors the source CSVs being re-parsed for a later cut-off week; several
# already-reported weeks get their observed/expected counts nudged as late
# registrations land, and a brand-new week (142 = "2022 week 36") is appended
# with its own observed/expected counts and %-change formulas.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing weekly observed/expected death counts (rows 120-141) ---
$ws.Range("AA120").Value = 273
$ws.Range("W123").Value = 448
$ws.Range("X125").Value = 585
$ws.Range("W126").Value = 455
$ws.Range("Z126").Value = 495
$ws.Range("AA130").Value = 229
$ws.Range("T132").Value = 52
$ws.Range("AA132").Value = 259
$ws.Range("V133").Value = 197
$ws.Range("X134").Value = 624
$ws.Range("W135").Value = 508
$ws.Range("AA135").Value = 245
$ws.Range("X136").Value = 596
$ws.Range("S137").Value = 213
$ws.Range("V137").Value = 197
$ws.Range("W137").Value = 452
$ws.Range("X137").Value = 611
$ws.Range("AA137").Value = 217
$ws.Range("Q138").Value = 119
$ws.Range("S138").Value = 211
$ws.Range("U138").Value = 394
$ws.Range("V138").Value = 175
$ws.Range("X138").Value = 590
$ws.Range("U139").Value = 379
$ws.Range("V139").Value = 209
$ws.Range("W139").Value = 497
$ws.Range("X139").Value = 600
$ws.Range("Z139").Value = 472
$ws.Range("AA139").Value = 203
$ws.Range("T140").Value = 40
$ws.Range("U140").Value = 355
$ws.Range("W140").Value = 436
$ws.Range("X140").Value = 591
$ws.Range("Z140").Value = 486
$ws.Range("AA140").Value = 230
$ws.Range("P141").Value = 88
$ws.Range("Q141").Value = 121
$ws.Range("R141").Value = 116
$ws.Range("S141").Value = 202
$ws.Range("T141").Value = 50
$ws.Range("U141").Value = 387
$ws.Range("V141").Value = 187
$ws.Range("W141").Value = 391
$ws.Range("X141").Value = 588
$ws.Range("Y141").Value = 71
$ws.Range("Z141").Value = 471
$ws.Range("AA141").Value = 238

# --- Append new week row 142 (2022 week 36) with observed + expected counts ---
$ws.Range("N142").Value = 2022
$ws.Range("O142").Value = 36
$ws.Range("P142").Value = 124
$ws.Range("Q142").Value = 115
$ws.Range("R142").Value = 89
$ws.Range("S142").Value = 197
$ws.Range("T142").Value = 47
$ws.Range("U142").Value = 368
$ws.Range("V142").Value = 183
$ws.Range("W142").Value = 484
$ws.Range("X142").Value = 594
$ws.Range("Y142").Value = 73
$ws.Range("Z142").Value = 442
$ws.Range("AA142").Value = 220
$ws.Range("AC142").Value = 2022
$ws.Range("AD142").Value = 36

# --- Add percentage-change formulas for row 142 (mirrors formulas in rows above) ---
$ws.Range("AE142").Formula = "=ROUND((P142-B142)/B142*100,2)"
$ws.Range("AF142").Formula = "=ROUND((Q142-C142)/C142*100,2)"
$ws.Range("AG142").Formula = "=ROUND((R142-D142)/D142*100,2)"
$ws.Range("AH142").Formula = "=ROUND((S142-E142)/E142*100,2)"
$ws.Range("AI142").Formula = "=ROUND((T142-F142)/F142*100,2)"
$ws.Range("AJ142").Formula = "=ROUND((U142-G142)/G142*100,2)"
$ws.Range("AK142").Formula = "=ROUND((V142-H142)/H142*100,2)"
$ws.Range("AL142").Formula = "=ROUND((W142-I142)/I142*100,2)"
$ws.Range("AM142").Formula = "=ROUND((X142-J142)/J142*100,2)"
$ws.Range("AN142").Formula = "=ROUND((Y142-K142)/K142*100,2)"
$ws.Range("AO142").Formula = "=ROUND((Z142-L142)/L142*100,2)"
$ws.Range("AP142").Formula = "=ROUND((AA142-M142)/M142*100,2)"

# --- Restore the view: scrolled down a bit further and the cursor left on
#     AH141 (where AH column's % figure was being double-checked) ---
$ws.Range("AH141").Select() | Out-Null
